$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels for B1 and E1 (kitchens_1 <-> living_rooms_2)
$ws.Range("B1").Value = "living_rooms_2"
$ws.Range("E1").Value = "kitchens_1"

# Update the one-hot block-order data (rows 2-7) to match the new column order
$values = @(
    @(0,0,0,1,0,0),
    @(0,0,0,0,1,0),
    @(1,0,0,0,0,0),
    @(0,0,1,0,0,0),
    @(0,0,0,0,0,1),
    @(0,1,0,0,0,0)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}

$wb.Save()
